$d = $word.ActiveDocument

$replacements = @(
    @("492÷9=54, 6", "925÷7=132, 1"),
    @("493÷7=70, 3", "658÷9=73, 1"),
    @("791÷9=87, 8", "706÷9=78, 4"),
    @("815÷8=101, 7", "360÷9=40, 0"),
    @("230÷6=38, 2", "144÷2=72, 0"),
    @("463÷9=51, 4", "880÷9=97, 7"),
    @("729÷8=91, 1", "344÷9=38, 2"),
    @("966÷4=241, 2", "478÷6=79, 4"),
    @("120÷4=30, 0", "826÷8=103, 2"),
    @("383÷9=42, 5", "975÷9=108, 3"),
    @("582÷9=64, 6", "662÷6=110, 2"),
    @("494÷9=54, 8", "276÷2=138, 0"),
    @("266÷6=44, 2", "112÷3=37, 1"),
    @("718÷9=79, 7", "682÷4=170, 2"),
    @("759÷3=253, 0", "400÷2=200, 0"),
    @("562÷3=187, 1", "217÷3=72, 1"),
    @("977÷7=139, 4", "720÷2=360, 0"),
    @("238÷2=119, 0", "471÷5=94, 1"),
    @("769÷4=192, 1", "413÷8=51, 5"),
    @("865÷3=288, 1", "785÷4=196, 1"),
    @("275÷4=68, 3", "741÷5=148, 1"),
    @("650÷2=325, 0", "119÷5=23, 4"),
    @("572÷6=95, 2", "641÷9=71, 2"),
    @("362÷8=45, 2", "400÷2=200, 0"),
    @("173÷7=24, 5", "994÷5=198, 4"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done."
